# "AGREGANDO CAMBIOS AL REPOSITORIO" - centre across C2:K2, merged, bold/italic/underline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:K2")
$range.HorizontalAlignment = -4108  # xlCenter
$range.Merge()

$title = $ws.Range("C2")
$title.Value = "AGREGANDO CAMBIOS AL REPOSITORIO "
$title.Font.Bold = $true
$title.Font.Italic = $true
$title.Font.Underline = 2  # xlUnderlineStyleSingle

$ws.PageSetup.PaperSize = 9       # xlPaperA4 (OOXML paperSize code 9)
$ws.PageSetup.Orientation = 1     # xlPortrait

$ws.Range("H6").Select() | Out-Null
